$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Km initiali
$ws.Range("B12").Value = 251008

# Data table (rows 14-43): Km_parcursi (B), Locul deplasarii (C), Observatii utilizator (D)
$data = @(
    @(14, 0,   "",               ""),
    @(15, 30,  "Acasa-Birou",    "Interes Serviciu"),
    @(16, 101, "Cluj-Dej",       "Interes Serviciu"),
    @(17, 0,   "",               ""),
    @(18, 0,   "",               ""),
    @(19, 356, "Cluj-Baia-Mare", "Interes Serviciu"),
    @(20, 152, "Cluj-Cmp. Turzii","Interes Serviciu"),
    @(21, 30,  "Acasa-Birou",    " "),
    @(22, 356, "Cluj-Baia-Mare", "Interes Serviciu"),
    @(23, 121, "Cluj-Turda",     "Interes Serviciu"),
    @(24, 0,   "",               ""),
    @(25, 0,   "",               ""),
    @(26, 0,   "",               ""),
    @(27, 85,  "Cluj-Apahida",   "Interes Serviciu"),
    @(28, 30,  "Acasa-Birou",    " "),
    @(29, 47,  "Cluj-Cluj",      "Interes Serviciu"),
    @(30, 30,  "Acasa-Birou",    " "),
    @(31, 0,   "",               ""),
    @(32, 0,   "",               ""),
    @(33, 85,  "Cluj-Apahida",   "Interes Serviciu"),
    @(34, 121, "Cluj-Turda",     "Interes Serviciu"),
    @(35, 30,  "Acasa-Birou",    " "),
    @(36, 257, "Cluj-Bistrita",  "Interes Serviciu"),
    @(37, 30,  "Acasa-Birou",    " "),
    @(38, 0,   "",               ""),
    @(39, 0,   "",               ""),
    @(40, 421, "Cluj-Satu-Mare", "Interes Serviciu"),
    @(41, 30,  "Acasa-Birou",    " "),
    @(42, 152, "Cluj-Cmp. Turzii","Interes Serviciu"),
    @(43, 85,  "Cluj-Apahida",   "Interes Serviciu")
)

foreach ($row in $data) {
    $r = $row[0]
    $km = $row[1]
    $loc = $row[2]
    $obs = $row[3]

    $ws.Cells.Item($r, 2).Value = $km
    if ($loc -ne "") {
        $ws.Cells.Item($r, 3).Value = $loc
    } else {
        $ws.Cells.Item($r, 3).Value = ""
    }
    if ($obs -ne "") {
        $ws.Cells.Item($r, 4).Value = $obs
    } else {
        $ws.Cells.Item($r, 4).Value = ""
    }
}

# Totals
$ws.Range("B44").Value = 2549
$ws.Range("B45").Value = 253557
